$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 9843.75
$ws.Range("I2").Value = 1763.6666
$ws.Range("K2").Value = 1763.6666
$ws.Range("M2").Value = -1650.6666
$ws.Range("H4").Value = 2877.5
$ws.Range("I4").Value = 3190.5557
$ws.Range("K4").Value = 3190.5557
$ws.Range("M4").Value = -3076.5557
$ws.Range("H9").Value = 1714235.5
$ws.Range("I9").Value = 135
$ws.Range("J9").Value = 2999811
$ws.Range("K9").Value = 135
$ws.Range("L9").Value = 2999811
$ws.Range("M9").Value = 34
$ws.Range("N9").Value = -3000149
$ws.Range("H28").Value = 2444.5293
$ws.Range("J28").Value = 680.1111
$ws.Range("L28").Value = 680.1111
$ws.Range("N28").Value = -1650.1111
$ws.Range("H103").Value = 1684
$ws.Range("J103").Value = 1447
$ws.Range("L103").Value = 4341
$ws.Range("N103").Value = -5513
$ws.Range("H138").Value = 2429.9583
$ws.Range("J138").Value = 2925.8333
$ws.Range("L138").Value = 8777.499899999999
$ws.Range("N138").Value = -19057.4999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 12000
$ws.Range("J23").Value = 12000
$ws.Range("L23").Value = 12000
$ws.Range("N23").Value = -12518
$ws.Range("H32").Value = 22050.309
$ws.Range("I32").Value = 22734.22
$ws.Range("J32").Value = 4952.5
$ws.Range("K32").Value = 22734.22
$ws.Range("L32").Value = 4952.5
$ws.Range("M32").Value = -22447.22
$ws.Range("N32").Value = -5526.5
$ws.Range("H50").Value = 386.33334
$ws.Range("I50").Value = 402
$ws.Range("J50").Value = 373.8
$ws.Range("K50").Value = 402
$ws.Range("L50").Value = 373.8
$ws.Range("M50").Value = 312
$ws.Range("N50").Value = -1801.8
$ws.Range("H61").Value = 6431.087
$ws.Range("I61").Value = 5487.6313
$ws.Range("J61").Value = 10912.5
$ws.Range("K61").Value = 5487.6313
$ws.Range("L61").Value = 10912.5
$ws.Range("M61").Value = -5275.6313
$ws.Range("N61").Value = -11336.5
$ws.Range("H63").Value = 9308.174000000001
$ws.Range("I63").Value = 2144
$ws.Range("K63").Value = 2144
$ws.Range("M63").Value = -1458
$ws.Range("H66").Value = 9308.174000000001
$ws.Range("I66").Value = 2144
$ws.Range("K66").Value = 10720
$ws.Range("M66").Value = -7288
$ws.Range("H136").Value = 6431.087
$ws.Range("I136").Value = 5487.6313
$ws.Range("J136").Value = 10912.5
$ws.Range("K136").Value = 16462.8939
$ws.Range("L136").Value = 32737.5
$ws.Range("M136").Value = -13912.8939
$ws.Range("N136").Value = -37837.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2169
$ws.Range("I20").Value = 2042.238
$ws.Range("K20").Value = 2042.238
$ws.Range("M20").Value = -1795.238
$ws.Range("H25").Value = 6007.0586
$ws.Range("I25").Value = 198.66667
$ws.Range("K25").Value = 198.66667
$ws.Range("M25").Value = 36.33332999999999
$ws.Range("H93").Value = 51947.5
$ws.Range("J93").Value = 51947.5
$ws.Range("L93").Value = 51947.5
$ws.Range("N93").Value = -55691.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34487050
$ws.Range("I31").Value = 111112350
$ws.Range("K31").Value = 111112350
$ws.Range("M31").Value = -111112055
$ws.Range("H34").Value = 34487050
$ws.Range("I34").Value = 111112350
$ws.Range("K34").Value = 111112350
$ws.Range("M34").Value = -111112148
$ws.Range("H97").Value = 59942.668
$ws.Range("J97").Value = 59942.668
$ws.Range("L97").Value = 59942.668
$ws.Range("N97").Value = -61924.668
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2720.7144
$ws.Range("I39").Value = 1507.6111
$ws.Range("K39").Value = 4522.8333
$ws.Range("M39").Value = -4228.8333
$ws.Range("H55").Value = 721.6
$ws.Range("I55").Value = 487.25
$ws.Range("J55").Value = 877.8333
$ws.Range("K55").Value = 1461.75
$ws.Range("L55").Value = 2633.4999
$ws.Range("M55").Value = -1284.75
$ws.Range("N55").Value = -2987.4999
$ws.Range("H68").Value = 403960
$ws.Range("I68").Value = 2000000
$ws.Range("J68").Value = 4950
$ws.Range("K68").Value = 6000000
$ws.Range("L68").Value = 14850
$ws.Range("M68").Value = -5999189
$ws.Range("N68").Value = -16472
$ws.Range("H71").Value = 403960
$ws.Range("I71").Value = 2000000
$ws.Range("J71").Value = 4950
$ws.Range("K71").Value = 18000000
$ws.Range("L71").Value = 44550
$ws.Range("M71").Value = -17995944
$ws.Range("N71").Value = -52662
$ws.Range("H87").Value = 2000
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 2000
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H103").Value = 4227
$ws.Range("J103").Value = 3665.6667
$ws.Range("L103").Value = 10997.0001
$ws.Range("N103").Value = -12755.0001
$ws.Range("H131").Value = 15877747
$ws.Range("J131").Value = 6218.1333
$ws.Range("L131").Value = 18654.3999
$ws.Range("N131").Value = -28734.3999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 1000
$ws.Range("K5").Value = 1000
$ws.Range("M5").Value = -888
$ws.Range("H31").Value = 2100
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H37").Value = 2100
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H70").Value = 13009
$ws.Range("I70").Value = 11664.3125
$ws.Range("K70").Value = 11664.3125
$ws.Range("M70").Value = -11394.3125
$ws.Range("H73").Value = 13009
$ws.Range("I73").Value = 11664.3125
$ws.Range("K73").Value = 11664.3125
$ws.Range("M73").Value = -10728.3125
$ws.Range("H99").Value = 9242.5
$ws.Range("I99").Value = 9242.5
$ws.Range("K99").Value = 9242.5
$ws.Range("M99").Value = -6996.5
$ws.Range("H140").Value = 110000
$ws.Range("J140").Value = 110000
$ws.Range("L140").Value = 110000
$ws.Range("N140").Value = -120360
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("H80").Value = 64111
$ws.Range("J80").Value = 64111
$ws.Range("L80").Value = 64111
$ws.Range("N80").Value = -66357
$ws.Range("H83").Value = 64111
$ws.Range("J83").Value = 64111
$ws.Range("L83").Value = 192333
$ws.Range("N83").Value = -203565
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 31342.5
$ws.Range("I70").Value = 31342.5
$ws.Range("K70").Value = 31342.5
$ws.Range("M70").Value = -31027.5
$ws.Range("H73").Value = 31342.5
$ws.Range("I73").Value = 31342.5
$ws.Range("K73").Value = 31342.5
$ws.Range("M73").Value = -30250.5
$ws.Range("H86").Value = 59125.668
$ws.Range("J86").Value = 59125.668
$ws.Range("L86").Value = 59125.668
$ws.Range("N86").Value = -61371.668
$ws.Range("H89").Value = 59125.668
$ws.Range("J89").Value = 59125.668
$ws.Range("L89").Value = 295628.34
$ws.Range("N89").Value = -306860.34
